$d = $word.ActiveDocument

# Paragraph 1 is the heading, currently split as:
#   [proofErr gramStart] "2" [proofErr gramEnd] " Objetivo do projeto"
# Target: a single run "2 Objetivo do projeto" with no proofErr markers
# (the grammar-check false positive on the lone "2" was cleared and the
# heading text merged into one run).

$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range

# Remember the run-level character formatting used by the heading text
# (Arial, bold, 12pt) before we touch anything.
$fontName = $r.Font.Name
$fontSize = $r.Font.Size
$fontBold = $r.Font.Bold

# Remove the whole paragraph range. This drops the text runs AND the
# gramStart/gramEnd proofErr markers that a plain Find/Replace would
# leave behind, while the freshly inserted paragraph mark keeps the
# same paragraph-level formatting (w:pPr/w:rPr).
$r.Delete()
$d.Range(0, 0).InsertParagraphBefore()

$newRange = $d.Range(0, 0)
$newRange.Font.Name = $fontName
$newRange.Font.Size = $fontSize
$newRange.Font.Bold = $fontBold
$newRange.InsertBefore("2 Objetivo do projeto")

Write-Host "Heading fixed."
